$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H7").Value = 1.2
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$s2 = $sc.Item(2)
$f = $s2.Formula
Write-Host "Formula: $f"
$s2.Formula = $f
